$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Bsg"
$ws.Cells.Item(2, 3).Value = "Slc16a7"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 38.60286033333333
$ws.Cells.Item(2, 8).Value = 115.808581
$ws.Cells.Item(2, 9).Value = 0.2650212684862838
$ws.Cells.Item(2, 10).Value = 0.2650212684862838
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 1.934841
$ws.Cells.Item(2, 14).Value = 5.804523
$ws.Cells.Item(2, 15).Value = 0.1811319458224623
$ws.Cells.Item(2, 16).Value = 0.1811319458224623
$ws.Cells.Item(2, 17).Value = 74.69039689020698
$ws.Cells.Item(2, 18).Value = 672.213572011863
$ws.Cells.Item(2, 19).Value = 0.04800381804525779
$ws.Cells.Item(2, 20).Value = 0.04800381804525779

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Bsg"
$ws.Cells.Item(3, 3).Value = "Slc16a7"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 38.60286033333333
$ws.Cells.Item(3, 8).Value = 115.808581
$ws.Cells.Item(3, 9).Value = 0.2650212684862838
$ws.Cells.Item(3, 10).Value = 0.2650212684862838
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.286000666666667
$ws.Cells.Item(3, 14).Value = 6.858002
$ws.Cells.Item(3, 15).Value = 0.2140060857221753
$ws.Cells.Item(3, 16).Value = 0.2140060857221753
$ws.Cells.Item(3, 17).Value = 88.24616445724023
$ws.Cells.Item(3, 18).Value = 794.215480115162
$ws.Cells.Item(3, 19).Value = 0.0567161643018753
$ws.Cells.Item(3, 20).Value = 0.05671616430187529

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Bsg"
$ws.Cells.Item(4, 3).Value = "Slc16a7"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 38.60286033333333
$ws.Cells.Item(4, 8).Value = 115.808581
$ws.Cells.Item(4, 9).Value = 0.2650212684862838
$ws.Cells.Item(4, 10).Value = 0.2650212684862838
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 6.461100666666667
$ws.Cells.Item(4, 14).Value = 19.383302
$ws.Cells.Item(4, 15).Value = 0.6048619684553623
$ws.Cells.Item(4, 16).Value = 0.6048619684553623
$ws.Cells.Item(4, 17).Value = 249.4169666349402
$ws.Cells.Item(4, 18).Value = 2244.752699714462
$ws.Cells.Item(4, 19).Value = 0.1603012861391507
$ws.Cells.Item(4, 20).Value = 0.1603012861391507

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Bsg"
$ws.Cells.Item(5, 3).Value = "Slc16a7"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 74.45592499999999
$ws.Cells.Item(5, 8).Value = 223.367775
$ws.Cells.Item(5, 9).Value = 0.5111642898850374
$ws.Cells.Item(5, 10).Value = 0.5111642898850374
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 1.934841
$ws.Cells.Item(5, 14).Value = 5.804523
$ws.Cells.Item(5, 15).Value = 0.1811319458224623
$ws.Cells.Item(5, 16).Value = 0.1811319458224623
$ws.Cells.Item(5, 17).Value = 144.060376382925
$ws.Cells.Item(5, 18).Value = 1296.543387446325
$ws.Cells.Item(5, 19).Value = 0.09258818246183399
$ws.Cells.Item(5, 20).Value = 0.09258818246183399

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Bsg"
$ws.Cells.Item(6, 3).Value = "Slc16a7"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 74.45592499999999
$ws.Cells.Item(6, 8).Value = 223.367775
$ws.Cells.Item(6, 9).Value = 0.5111642898850374
$ws.Cells.Item(6, 10).Value = 0.5111642898850374
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.286000666666667
$ws.Cells.Item(6, 14).Value = 6.858002
$ws.Cells.Item(6, 15).Value = 0.2140060857221753
$ws.Cells.Item(6, 16).Value = 0.2140060857221753
$ws.Cells.Item(6, 17).Value = 170.2062941872833
$ws.Cells.Item(6, 18).Value = 1531.85664768555
$ws.Cells.Item(6, 19).Value = 0.1093922688392522
$ws.Cells.Item(6, 20).Value = 0.1093922688392522

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Bsg"
$ws.Cells.Item(7, 3).Value = "Slc16a7"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 74.45592499999999
$ws.Cells.Item(7, 8).Value = 223.367775
$ws.Cells.Item(7, 9).Value = 0.5111642898850374
$ws.Cells.Item(7, 10).Value = 0.5111642898850374
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 6.461100666666667
$ws.Cells.Item(7, 14).Value = 19.383302
$ws.Cells.Item(7, 15).Value = 0.6048619684553623
$ws.Cells.Item(7, 16).Value = 0.6048619684553623
$ws.Cells.Item(7, 17).Value = 481.0672266547833
$ws.Cells.Item(7, 18).Value = 4329.60503989305
$ws.Cells.Item(7, 19).Value = 0.3091838385839512
$ws.Cells.Item(7, 20).Value = 0.3091838385839512

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Bsg"
$ws.Cells.Item(8, 3).Value = "Slc16a7"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 32.60069533333333
$ws.Cells.Item(8, 8).Value = 97.80208599999999
$ws.Cells.Item(8, 9).Value = 0.2238144416286788
$ws.Cells.Item(8, 10).Value = 0.2238144416286788
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 1.934841
$ws.Cells.Item(8, 14).Value = 5.804523
$ws.Cells.Item(8, 15).Value = 0.1811319458224623
$ws.Cells.Item(8, 16).Value = 0.1811319458224623
$ws.Cells.Item(8, 17).Value = 63.07716195944198
$ws.Cells.Item(8, 18).Value = 567.6944576349779
$ws.Cells.Item(8, 19).Value = 0.04053994531537049
$ws.Cells.Item(8, 20).Value = 0.04053994531537049

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Bsg"
$ws.Cells.Item(9, 3).Value = "Slc16a7"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 32.60069533333333
$ws.Cells.Item(9, 8).Value = 97.80208599999999
$ws.Cells.Item(9, 9).Value = 0.2238144416286788
$ws.Cells.Item(9, 10).Value = 0.2238144416286788
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.286000666666667
$ws.Cells.Item(9, 14).Value = 6.858002
$ws.Cells.Item(9, 15).Value = 0.2140060857221753
$ws.Cells.Item(9, 16).Value = 0.2140060857221753
$ws.Cells.Item(9, 17).Value = 74.52521126579688
$ws.Cells.Item(9, 18).Value = 670.726901392172
$ws.Cells.Item(9, 19).Value = 0.04789765258104783
$ws.Cells.Item(9, 20).Value = 0.04789765258104783

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Bsg"
$ws.Cells.Item(10, 3).Value = "Slc16a7"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 32.60069533333333
$ws.Cells.Item(10, 8).Value = 97.80208599999999
$ws.Cells.Item(10, 9).Value = 0.2238144416286788
$ws.Cells.Item(10, 10).Value = 0.2238144416286788
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 6.461100666666667
$ws.Cells.Item(10, 14).Value = 19.383302
$ws.Cells.Item(10, 15).Value = 0.6048619684553623
$ws.Cells.Item(10, 16).Value = 0.6048619684553623
$ws.Cells.Item(10, 17).Value = 210.6363743519969
$ws.Cells.Item(10, 18).Value = 1895.727369167972
$ws.Cells.Item(10, 19).Value = 0.1353768437322604
$ws.Cells.Item(10, 20).Value = 0.1353768437322604

